$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they stay as plain text
# (matching the original inline-string storage of these columns).
$textCells = @("D5", "D6", "D8", "D10", "D11", "D14", "D15", "D20", "D22", "D25", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D43", "D45", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.104.48"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "1.931.49"
$ws.Range("E3").Value = "  -5.03%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "240.52"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").Value = "0.604"
$ws.Range("E6").Value = "  -5.42%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "55.66"
$ws.Range("E8").Value = "  -11.27%  "
$ws.Range("E9").Value = "  -8.97%  "
$ws.Range("D10").Value = "55.26"
$ws.Range("E10").Value = "  -4.85%  "
$ws.Range("D11").Value = "0.0822"
$ws.Range("E11").Value = "  +4.60%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.217.81"
$ws.Range("E13").Value = "  -4.83%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.811"
$ws.Range("E14").Value = "  -9.55%  "
$ws.Range("D15").Value = "20.76"
$ws.Range("E15").Value = "  -11.61%  "
$ws.Range("E16").Value = "  -8.96%  "
$ws.Range("E17").Value = "  -7.70%  "
$ws.Range("D18").Value = "1.914.61"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("D19").Value = "36.030.14"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "69.26"
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").Value = "224.76"
$ws.Range("E22").Value = "  -5.07%  "
$ws.Range("E23").Value = "  -8.81%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  -5.13%  "
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  -6.85%  "
$ws.Range("D28").Value = "162.37"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "19.11"
$ws.Range("E29").Value = "  -6.22%  "
$ws.Range("D30").Value = "0.116"
$ws.Range("E30").Value = "  -17.76%  "
$ws.Range("D31").Value = "0.117"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("E32").Value = "  -6.09%  "
$ws.Range("D33").Value = "4.64"
$ws.Range("E33").Value = "  -8.22%  "
$ws.Range("D34").Value = "0.0618"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  -7.09%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.79"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "5.86"
$ws.Range("E38").Value = "  -11.44%  "
$ws.Range("D39").Value = "2.13"
$ws.Range("E39").Value = "  -10.79%  "
$ws.Range("E40").Value = "  -11.36%  "
$ws.Range("D41").Value = "0.0965"
$ws.Range("E41").Value = "  -4.91%  "
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "1.15"
$ws.Range("E43").Value = "  -9.09%  "
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").Value = "15.40"
$ws.Range("E45").Value = "  -9.28%  "
$ws.Range("D46").Value = "1.333.69"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "87.01"
$ws.Range("E48").Value = "  -7.68%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.15"
$ws.Range("E49").Value = "  -7.39%  "
$ws.Range("D50").Value = "2.80"
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("D51").Value = "45.75"
$ws.Range("E51").Value = "  +0.83%  "
